$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text or non-ambiguous (safe to assign directly)
$plainUpdates = @{
    "D2" = "29.480.73"
    "D3" = "1.849.51"
    "E3" = "  -0.65%  "
    "E4" = "  -0.04%  "
    "E5" = "  -1.28%  "
    "E6" = "  -2.38%  "
    "E7" = "  +0.00%  "
    "E8" = "  +0.41%  "
    "E9" = "  -0.31%  "
    "E10" = "  -0.08%  "
    "E11" = "  -1.27%  "
    "E12" = "  +0.11%  "
    "D13" = "1.908.84"
    "E13" = "  +2.54%  "
    "E14" = "  -0.79%  "
    "E15" = "  -1.14%  "
    "E16" = "  -0.25%  "
    "E17" = "  -1.56%  "
    "D18" = "2.156.58"
    "E18" = "  +2.21%  "
    "E19" = "  +1.85%  "
    "D20" = "29.560.42"
    "E20" = "  -0.68%  "
    "E21" = "  -0.77%  "
    "E22" = "  -1.65%  "
    "E23" = "  +0.00%  "
    "E24" = "  +1.42%  "
    "E25" = "  -0.04%  "
    "E26" = "  -2.34%  "
    "E27" = "  -2.31%  "
    "E28" = "  -1.52%  "
    "E29" = "  -1.26%  "
    "E30" = "  -0.67%  "
    "E31" = "  -5.42%  "
    "E32" = "  -2.35%  "
    "E33" = "  -1.45%  "
    "E34" = "  -2.03%  "
    "E35" = "  -0.76%  "
    "E36" = "  -0.39%  "
    "E37" = "  -1.42%  "
    "E38" = "  -0.56%  "
    "B39" = "MXToken"
    "C39" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "E39" = "  -1.12%  "
    "B40" = "Maker"
    "C40" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D40" = "1.239.35"
    "E40" = "  +2.72%  "
    "E41" = "  -0.40%  "
    "E42" = "  -0.93%  "
    "E43" = "  -2.65%  "
    "D44" = "2.065.43"
    "E44" = "  +2.20%  "
    "E45" = "  -0.05%  "
    "E46" = "  -0.29%  "
    "E47" = "  +1.20%  "
    "E48" = "  +8.34%  "
    "E49" = "  -0.18%  "
    "E50" = "  -1.30%  "
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Cells whose new values look like numbers but must remain text (matches source format),
# so force text formatting, assign, then restore the original style to avoid leaving
# a stray number-format override on the cell.
$textForcedUpdates = @{
    "D4" = "0.9991"
    "D5" = "241.95"
    "D6" = "0.6261"
    "D8" = "47.96"
    "D9" = "0.07529"
    "D10" = "0.2973"
    "D11" = "24.27"
    "D12" = "0.07685"
    "D14" = "5.006"
    "D15" = "0.6856"
    "D16" = "83.86"
    "D17" = "0.000009733"
    "D19" = "6.232"
    "D21" = "234.36"
    "D22" = "12.47"
    "D23" = "1.0000"
    "D24" = "7.614"
    "D26" = "155.68"
    "D27" = "0.1391"
    "D28" = "8.427"
    "D29" = "17.72"
    "D30" = "1.481"
    "D31" = "0.05840"
    "D33" = "4.104"
    "D34" = "4.020"
    "D35" = "1.882"
    "D37" = "0.7207"
    "D38" = "2.589"
    "D39" = "2.796"
    "D41" = "0.01779"
    "D42" = "0.9126"
    "D43" = "6.119"
    "D45" = "0.9997"
    "D46" = "101.72"
    "D48" = "7.274"
    "D49" = "9.163"
    "D51" = "0.4028"
}

foreach ($ref in $textForcedUpdates.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$ref]
    $cell.Style = $origStyle
}
